# Update "想去人数" (interested-count) values that changed in the
# upstream data refresh (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 12
$ws1.Range("F6").Value  = 57
$ws1.Range("F8").Value  = 11437
$ws1.Range("F9").Value  = 4327
$ws1.Range("F11").Value = 32
$ws1.Range("F12").Value = 18
$ws1.Range("F16").Value = 27
$ws1.Range("F20").Value = 11173

# --- Sheet "全部类型" ---------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 12
$ws4.Range("F6").Value  = 57
$ws4.Range("F8").Value  = 11437
$ws4.Range("F9").Value  = 4327
$ws4.Range("F11").Value = 32
$ws4.Range("F12").Value = 18
$ws4.Range("F17").Value = 27
$ws4.Range("F21").Value = 11173
